$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.805999999999999
$ws.Range("A9").Value = -21.847
$ws.Range("D9").Value = -8.210000000000001
$ws.Range("A18").Value = -22.095
$ws.Range("A20").Value = -20.519
$ws.Range("D23").Value = -7.992999999999999
$ws.Range("D24").Value = -6.946000000000001
$ws.Range("D26").Value = -7.48
$ws.Range("A27").Value = -21.942
$ws.Range("D34").Value = -7.730999999999999
$ws.Range("D35").Value = -7.877999999999998
$ws.Range("D48").Value = -7.892999999999999
$ws.Range("D52").Value = -7.953
$ws.Range("D66").Value = -7.486
$ws.Range("D67").Value = -7.616
$ws.Range("A69").Value = -21.52
$ws.Range("A76").Value = -20.043
$ws.Range("D80").Value = -7.906000000000001
$ws.Range("A82").Value = -22.067
$ws.Range("D99").Value = -8.253
